$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Trimestre (column C) dates for all data rows (force text, not a date serial)
$ws.Range("C2:C10").Value = "'01/01/2024"
$ws.Range("C2:C10").Style = "Normal"

# Update Valor (column D) values
$ws.Range("D2").Value = 56.34597304920088
$ws.Range("D3").Value = 55.75444179262795
$ws.Range("D4").Value = 55.01788344203661
$ws.Range("D5").Value = 54.78366030636926
$ws.Range("D6").Value = 54.20309929876201
$ws.Range("D7").Value = 53.84615384615385
$ws.Range("D8").Value = 46.585672392124
$ws.Range("D9").Value = 43.09459621889864
$ws.Range("D10").Value = 50.2716236441915

# Update Colocação (column E) rank for Sergipe
$ws.Range("E8").Value = "15º"
